$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 489.27274
$ws.Range("I5").Value = 489.27274
$ws.Range("K5").Value = 489.27274
$ws.Range("M5").Value = -374.27274
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 1894920.4
$ws.Range("I12").Value = 2525894
$ws.Range("J12").Value = 1999.6666
$ws.Range("K12").Value = 2525894
$ws.Range("L12").Value = 1999.6666
$ws.Range("M12").Value = -2525724
$ws.Range("N12").Value = -2339.6666
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H18").Value = 2874.111
$ws.Range("I18").Value = 2874.111
$ws.Range("K18").Value = 2874.111
$ws.Range("M18").Value = -2590.111
$ws.Range("H58").Value = 1587
$ws.Range("I58").Value = 441.33334
$ws.Range("K58").Value = 1324.00002
$ws.Range("M58").Value = -1174.00002
$ws.Range("H92").Value = 2512.4
$ws.Range("I92").Value = 1895.8572
$ws.Range("K92").Value = 1895.8572
$ws.Range("M92").Value = -647.8571999999999
$ws.Range("H111").Value = 10103091
$ws.Range("I111").Value = 11113150
$ws.Range("K111").Value = 33339450
$ws.Range("M111").Value = -33336383
$ws.Range("H118").Value = 90909750
$ws.Range("I118").Value = 166667230
$ws.Range("J118").Value = 785
$ws.Range("K118").Value = 500001690
$ws.Range("L118").Value = 2355
$ws.Range("M118").Value = -500000033
$ws.Range("N118").Value = -5669
$ws.Range("H137").Value = 95797.89999999999
$ws.Range("I137").Value = 199815.33
$ws.Range("K137").Value = 599445.99
$ws.Range("M137").Value = -596895.99
$ws.Range("H138").Value = 2576.5088
$ws.Range("I138").Value = 1304.8334
$ws.Range("J138").Value = 3501.3635
$ws.Range("K138").Value = 3914.5002
$ws.Range("L138").Value = 10504.0905
$ws.Range("M138").Value = 1225.4998
$ws.Range("N138").Value = -20784.0905
$ws.Range("H141").Value = 36294.8
$ws.Range("I141").Value = 27190
$ws.Range("K141").Value = 81570
$ws.Range("M141").Value = -76390

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1664350.1
$ws.Range("I2").Value = 2175758.2
$ws.Range("J2").Value = 2273.625
$ws.Range("K2").Value = 2175758.2
$ws.Range("L2").Value = 2273.625
$ws.Range("M2").Value = -2175645.2
$ws.Range("N2").Value = -2499.625
$ws.Range("H32").Value = 11774.966
$ws.Range("I32").Value = 7313.0234
$ws.Range("J32").Value = 24565.867
$ws.Range("K32").Value = 7313.0234
$ws.Range("L32").Value = 24565.867
$ws.Range("M32").Value = -7026.0234
$ws.Range("N32").Value = -25139.867
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H61").Value = 4944.9
$ws.Range("I61").Value = 4935.7646
$ws.Range("J61").Value = 4996.6665
$ws.Range("K61").Value = 4935.7646
$ws.Range("L61").Value = 4996.6665
$ws.Range("M61").Value = -4723.7646
$ws.Range("N61").Value = -5420.6665
$ws.Range("H63").Value = 5670.5293
$ws.Range("I63").Value = 3280
$ws.Range("J63").Value = 9085.571
$ws.Range("K63").Value = 3280
$ws.Range("L63").Value = 9085.571
$ws.Range("M63").Value = -2594
$ws.Range("N63").Value = -10457.571
$ws.Range("H66").Value = 5670.5293
$ws.Range("I66").Value = 3280
$ws.Range("J66").Value = 9085.571
$ws.Range("K66").Value = 16400
$ws.Range("L66").Value = 45427.855
$ws.Range("M66").Value = -12968
$ws.Range("N66").Value = -52291.855
$ws.Range("H74").Value = 102906.39
$ws.Range("I74").Value = 91748.5
$ws.Range("K74").Value = 91748.5
$ws.Range("M74").Value = -90874.5
$ws.Range("H77").Value = 102906.39
$ws.Range("I77").Value = 91748.5
$ws.Range("K77").Value = 458742.5
$ws.Range("M77").Value = -454374.5
$ws.Range("H88").Value = 1779.5714
$ws.Range("J88").Value = 1248.25
$ws.Range("L88").Value = 1248.25
$ws.Range("N88").Value = -2060.25
$ws.Range("H91").Value = 1779.5714
$ws.Range("J91").Value = 1248.25
$ws.Range("L91").Value = 1248.25
$ws.Range("N91").Value = -4056.25
$ws.Range("H97").Value = 1702116.2
$ws.Range("I97").Value = 2155905.8
$ws.Range("K97").Value = 2155905.8
$ws.Range("M97").Value = -2155409.8
$ws.Range("H110").Value = 868939.75
$ws.Range("I110").Value = 1029594.44
$ws.Range("K110").Value = 1029594.44
$ws.Range("M110").Value = -1027549.44
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H116").Value = 1664350.1
$ws.Range("I116").Value = 2175758.2
$ws.Range("J116").Value = 2273.625
$ws.Range("K116").Value = 2175758.2
$ws.Range("L116").Value = 2273.625
$ws.Range("M116").Value = -2173464.2
$ws.Range("N116").Value = -6861.625
$ws.Range("H132").Value = 3114.5435
$ws.Range("I132").Value = 1654.6111
$ws.Range("J132").Value = 4053.0715
$ws.Range("K132").Value = 4963.8333
$ws.Range("L132").Value = 12159.2145
$ws.Range("M132").Value = -2433.8333
$ws.Range("N132").Value = -17219.2145
$ws.Range("H136").Value = 4944.9
$ws.Range("I136").Value = 4935.7646
$ws.Range("J136").Value = 4996.6665
$ws.Range("K136").Value = 14807.2938
$ws.Range("L136").Value = 14989.9995
$ws.Range("M136").Value = -12257.2938
$ws.Range("N136").Value = -20089.9995
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 50000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 50000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -44900
$ws.Range("N137").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1664350.1
$ws.Range("I3").Value = 2175758.2
$ws.Range("J3").Value = 2273.625
$ws.Range("K3").Value = 2175758.2
$ws.Range("L3").Value = 2273.625
$ws.Range("M3").Value = -2175644.2
$ws.Range("N3").Value = -2501.625
$ws.Range("H38").Value = 39000
$ws.Range("J38").Value = 39000
$ws.Range("L38").Value = 39000
$ws.Range("N38").Value = -39832
$ws.Range("H86").Value = 13981128
$ws.Range("I86").Value = 18057742
$ws.Range("J86").Value = 4167.2856
$ws.Range("K86").Value = 18057742
$ws.Range("L86").Value = 4167.2856
$ws.Range("M86").Value = -18056619
$ws.Range("N86").Value = -6413.2856
$ws.Range("H89").Value = 13981128
$ws.Range("I89").Value = 18057742
$ws.Range("J89").Value = 4167.2856
$ws.Range("K89").Value = 90288710
$ws.Range("L89").Value = 20836.428
$ws.Range("M89").Value = -90283094
$ws.Range("N89").Value = -32068.428
$ws.Range("H94").Value = 2945140.5
$ws.Range("I94").Value = 3126399.5
$ws.Range("J94").Value = 45000
$ws.Range("K94").Value = 3126399.5
$ws.Range("L94").Value = 45000
$ws.Range("M94").Value = -3125948.5
$ws.Range("N94").Value = -45902
$ws.Range("H99").Value = 8405283
$ws.Range("I99").Value = 12988185
$ws.Range("K99").Value = 12988185
$ws.Range("M99").Value = -12986687
$ws.Range("H105").Value = 4809091.5
$ws.Range("I105").Value = 5209599
$ws.Range("K105").Value = 5209599
$ws.Range("M105").Value = -5207852
$ws.Range("H107").Value = 2101763.8
$ws.Range("I107").Value = 2551833.8
$ws.Range("K107").Value = 2551833.8
$ws.Range("M107").Value = -2549913.8
$ws.Range("H134").Value = 9970
$ws.Range("I134").Value = 1856.1428
$ws.Range("K134").Value = 5568.428400000001
$ws.Range("M134").Value = -3033.428400000001
$ws.Range("H140").Value = 79999.5
$ws.Range("J140").Value = 79999.5
$ws.Range("L140").Value = 79999.5
$ws.Range("N140").Value = -90359.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21728.96
$ws.Range("I31").Value = 3251.5
$ws.Range("J31").Value = 27414.334
$ws.Range("K31").Value = 3251.5
$ws.Range("L31").Value = 27414.334
$ws.Range("M31").Value = -2956.5
$ws.Range("N31").Value = -28004.334
$ws.Range("H34").Value = 21728.96
$ws.Range("I34").Value = 3251.5
$ws.Range("J34").Value = 27414.334
$ws.Range("K34").Value = 3251.5
$ws.Range("L34").Value = 27414.334
$ws.Range("M34").Value = -3049.5
$ws.Range("N34").Value = -27818.334
$ws.Range("H35").Value = 5316.6665
$ws.Range("I35").Value = 2980.2
$ws.Range("K35").Value = 2980.2
$ws.Range("M35").Value = -2686.2
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2112
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 14999.5
$ws.Range("I38").Value = 14999.5
$ws.Range("K38").Value = 14999.5
$ws.Range("M38").Value = -14622.5
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2340
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 20625
$ws.Range("J41").Value = 35000
$ws.Range("L41").Value = 35000
$ws.Range("N41").Value = -35856
$ws.Range("H46").Value = 14999.5
$ws.Range("I46").Value = 14999.5
$ws.Range("K46").Value = 14999.5
$ws.Range("M46").Value = -14788.5
$ws.Range("H50").Value = 12000
$ws.Range("J50").Value = 12000
$ws.Range("L50").Value = 12000
$ws.Range("N50").Value = -13250
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 5170.6772
$ws.Range("I58").Value = 5930.2856
$ws.Range("J58").Value = 3575.5
$ws.Range("K58").Value = 5930.2856
$ws.Range("L58").Value = 3575.5
$ws.Range("M58").Value = -5727.2856
$ws.Range("N58").Value = -3981.5
$ws.Range("H60").Value = 13000
$ws.Range("I60").Value = 13000
$ws.Range("K60").Value = 13000
$ws.Range("M60").Value = -12489
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 997.5
$ws.Range("I62").Value = 997.5
$ws.Range("K62").Value = 997.5
$ws.Range("M62").Value = -373.5
$ws.Range("H65").Value = 997.5
$ws.Range("I65").Value = 997.5
$ws.Range("K65").Value = 4987.5
$ws.Range("M65").Value = -1867.5
$ws.Range("H94").Value = 1757.1666
$ws.Range("I94").Value = 750
$ws.Range("J94").Value = 2260.75
$ws.Range("K94").Value = 750
$ws.Range("L94").Value = 2260.75
$ws.Range("M94").Value = -299
$ws.Range("N94").Value = -3162.75
$ws.Range("H132").Value = 91839.95
$ws.Range("I132").Value = 60337.766
$ws.Range("J132").Value = 225724.25
$ws.Range("K132").Value = 181013.298
$ws.Range("L132").Value = 677172.75
$ws.Range("M132").Value = -178483.298
$ws.Range("N132").Value = -682232.75
$ws.Range("H134").Value = 30815.363
$ws.Range("I134").Value = 92211.3
$ws.Range("K134").Value = 276633.9
$ws.Range("M134").Value = -274098.9
$ws.Range("H136").Value = 5170.6772
$ws.Range("I136").Value = 5930.2856
$ws.Range("J136").Value = 3575.5
$ws.Range("K136").Value = 17790.8568
$ws.Range("L136").Value = 10726.5
$ws.Range("M136").Value = -15240.8568
$ws.Range("N136").Value = -15826.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 77895.16
$ws.Range("J5").Value = 334796
$ws.Range("L5").Value = 1004388
$ws.Range("N5").Value = -1004612
$ws.Range("H24").Value = 866.6667
$ws.Range("I24").Value = 840.2
$ws.Range("J24").Value = 999
$ws.Range("K24").Value = 2520.6
$ws.Range("L24").Value = 2997
$ws.Range("M24").Value = -2290.6
$ws.Range("N24").Value = -3457
$ws.Range("H68").Value = 1370.4667
$ws.Range("I68").Value = 686.7273
$ws.Range("K68").Value = 2060.1819
$ws.Range("M68").Value = -1249.1819
$ws.Range("H71").Value = 1370.4667
$ws.Range("I71").Value = 686.7273
$ws.Range("K71").Value = 6180.545700000001
$ws.Range("M71").Value = -2124.545700000001
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
$ws.Range("H114").Value = 866.9524
$ws.Range("I114").Value = 148
$ws.Range("J114").Value = 1154.5333
$ws.Range("K114").Value = 444
$ws.Range("L114").Value = 3463.5999
$ws.Range("M114").Value = 2810
$ws.Range("N114").Value = -9971.599900000001
$ws.Range("H117").Value = 3126.0908
$ws.Range("I117").Value = 4139
$ws.Range("J117").Value = 2282
$ws.Range("K117").Value = 12417
$ws.Range("L117").Value = 6846
$ws.Range("M117").Value = -8975
$ws.Range("N117").Value = -13730
$ws.Range("H118").Value = 1030.6
$ws.Range("I118").Value = 1030.6
$ws.Range("K118").Value = 3091.8
$ws.Range("M118").Value = -1848.8
$ws.Range("H131").Value = 13444826
$ws.Range("I131").Value = 4387334.5
$ws.Range("K131").Value = 13162003.5
$ws.Range("M131").Value = -13156963.5
$ws.Range("H132").Value = 1624.4445
$ws.Range("I132").Value = 1463.7142
$ws.Range("J132").Value = 1797.5385
$ws.Range("K132").Value = 13173.4278
$ws.Range("L132").Value = 16177.8465
$ws.Range("M132").Value = -10643.4278
$ws.Range("N132").Value = -21237.8465
$ws.Range("H135").Value = 77895.16
$ws.Range("J135").Value = 334796
$ws.Range("L135").Value = 3013164
$ws.Range("N135").Value = -3018234
$ws.Range("H137").Value = 3706.6
$ws.Range("I137").Value = 3556
$ws.Range("K137").Value = 10668
$ws.Range("M137").Value = -5568
$ws.Range("H141").Value = 2896.2856
$ws.Range("J141").Value = 3500
$ws.Range("L141").Value = 10500
$ws.Range("N141").Value = -20860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 19800
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 19800
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 19800
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -20380
$ws.Range("H80").Value = 27725428
$ws.Range("I80").Value = 37470080
$ws.Range("K80").Value = 37470080
$ws.Range("M80").Value = -37469082
$ws.Range("H83").Value = 27725428
$ws.Range("I83").Value = 37470080
$ws.Range("K83").Value = 187350400
$ws.Range("M83").Value = -187345408
$ws.Range("H99").Value = 16686.2
$ws.Range("I99").Value = 17060.25
$ws.Range("J99").Value = 15190
$ws.Range("K99").Value = 17060.25
$ws.Range("L99").Value = 15190
$ws.Range("M99").Value = -14814.25
$ws.Range("N99").Value = -19682
$ws.Range("H102").Value = 5868544.5
$ws.Range("I102").Value = 9260267
$ws.Range("J102").Value = 2168483.2
$ws.Range("K102").Value = 9260267
$ws.Range("L102").Value = 2168483.2
$ws.Range("M102").Value = -9258645
$ws.Range("N102").Value = -2171727.2
$ws.Range("H107").Value = 575.1111
$ws.Range("I107").Value = 797.6
$ws.Range("J107").Value = 489.53845
$ws.Range("K107").Value = 797.6
$ws.Range("L107").Value = 489.53845
$ws.Range("M107").Value = 1122.4
$ws.Range("N107").Value = -4329.53845
$ws.Range("H113").Value = 6669261
$ws.Range("I113").Value = 11112804
$ws.Range("J113").Value = 3946.6
$ws.Range("K113").Value = 11112804
$ws.Range("L113").Value = 3946.6
$ws.Range("M113").Value = -11110634
$ws.Range("N113").Value = -8286.6
$ws.Range("H122").Value = 529229.3
$ws.Range("I122").Value = 894740
$ws.Range("J122").Value = 7071.143
$ws.Range("K122").Value = 2684220
$ws.Range("L122").Value = 21213.429
$ws.Range("M122").Value = -2681770
$ws.Range("N122").Value = -26113.429
$ws.Range("H132").Value = 2520.532
$ws.Range("I132").Value = 2501.1516
$ws.Range("J132").Value = 2566.2144
$ws.Range("K132").Value = 7503.4548
$ws.Range("L132").Value = 7698.6432
$ws.Range("M132").Value = -4973.4548
$ws.Range("N132").Value = -12758.6432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 29900
$ws.Range("J6").Value = 29900
$ws.Range("L6").Value = 29900
$ws.Range("N6").Value = -30124
$ws.Range("H42").Value = 8300
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 8300
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 8300
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -9426
$ws.Range("H49").Value = 8300
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 8300
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 8300
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -8594
$ws.Range("H82").Value = 52779450
$ws.Range("I82").Value = 105557380
$ws.Range("J82").Value = 1520.8
$ws.Range("K82").Value = 105557380
$ws.Range("L82").Value = 1520.8
$ws.Range("M82").Value = -105557019
$ws.Range("N82").Value = -2242.8
$ws.Range("H85").Value = 52779450
$ws.Range("I85").Value = 105557380
$ws.Range("J85").Value = 1520.8
$ws.Range("K85").Value = 105557380
$ws.Range("L85").Value = 1520.8
$ws.Range("M85").Value = -105556132
$ws.Range("N85").Value = -4016.8
$ws.Range("H100").Value = 3699.5881
$ws.Range("I100").Value = 2628.4285
$ws.Range("K100").Value = 2628.4285
$ws.Range("M100").Value = -2087.4285
$ws.Range("H125").Value = 73997.5
$ws.Range("J125").Value = 73997.5
$ws.Range("L125").Value = 73997.5
$ws.Range("N125").Value = -83837.5
$ws.Range("H132").Value = 6323.9575
$ws.Range("I132").Value = 6277.864
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 18833.592
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -16303.592
$ws.Range("N132").Value = -26060
$ws.Range("H136").Value = 55063.977
$ws.Range("I136").Value = 128493.31
$ws.Range("K136").Value = 385479.93
$ws.Range("M136").Value = -382929.93

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 26000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 26000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 26000
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -26406
$ws.Range("H81").Value = 9264103
$ws.Range("I81").Value = 15154032
$ws.Range("J81").Value = 8499.857
$ws.Range("K81").Value = 30308064
$ws.Range("L81").Value = 16999.714
$ws.Range("M81").Value = -30307003
$ws.Range("N81").Value = -19121.714
$ws.Range("H84").Value = 9264103
$ws.Range("I84").Value = 15154032
$ws.Range("J84").Value = 8499.857
$ws.Range("K84").Value = 151540320
$ws.Range("L84").Value = 84998.57000000001
$ws.Range("M84").Value = -151535016
$ws.Range("N84").Value = -95606.57000000001
$ws.Range("H122").Value = 4906.3887
$ws.Range("I122").Value = 2586.25
$ws.Range("J122").Value = 9546.666999999999
$ws.Range("K122").Value = 7758.75
$ws.Range("L122").Value = 28640.001
$ws.Range("M122").Value = -5308.75
$ws.Range("N122").Value = -33540.001
$ws.Range("H125").Value = 79799
$ws.Range("J125").Value = 79799
$ws.Range("L125").Value = 79799
$ws.Range("N125").Value = -89639
$ws.Range("H126").Value = 2200.1667
$ws.Range("I126").Value = 2250.2
$ws.Range("J126").Value = 2137.625
$ws.Range("K126").Value = 6750.599999999999
$ws.Range("L126").Value = 6412.875
$ws.Range("M126").Value = -4280.599999999999
$ws.Range("N126").Value = -11352.875
$ws.Range("H132").Value = 18723498
$ws.Range("I132").Value = 23259584
$ws.Range("J132").Value = 991524.5600000001
$ws.Range("K132").Value = 69778752
$ws.Range("L132").Value = 2974573.68
$ws.Range("M132").Value = -69776222
$ws.Range("N132").Value = -2979633.68
